# The diff removes the explicit "Don't add space between paragraphs of
# the same style" direct-formatting override (<w:contextualSpacing w:val="0"/>)
# from every paragraph's <w:pPr> in the document. In the Word object model
# this is the ParagraphFormat.ContextualSpacing property, so the edit is:
# turn it off (it was explicitly "0"/False) on every paragraph, which tells
# Word it no longer needs to carry the direct-formatting override at all.

$d = $word.ActiveDocument

foreach ($para in $d.Paragraphs) {
    $para.Format.ContextualSpacing = $false
}

Write-Host "ContextualSpacing cleared on $($d.Paragraphs.Count) paragraphs"
